# Insert a new "Version" column at the very left of the currency table.
# Every existing column (Code, Description, Definition, Numeric_Code) shifts
# one column to the right (A->B, B->C, C->D, D->E). The new column A gets the
# header "Version" in row 1, and the constant value "4217:2015" for every
# data row below it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing columns A:D to B:E, inserting a blank column A.
$ws.Columns("A:A").Insert()

# Figure out how many rows of data now exist (header + data rows).
$lastRow = $ws.UsedRange.Rows.Count

# Header for the new column.
$ws.Range("A1").Value2 = "Version"

# Constant value for every data row in the new column.
if ($lastRow -ge 2) {
    $ws.Range("A2:A$lastRow").Value2 = "4217:2015"
}
